# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- Update the "time_taken" (column F) timestamps on the data sheet ---
$timestamps = @{
    2  = "2021-10-05 14:33:18.375812"
    3  = "2021-10-05 14:33:18.375820"
    4  = "2021-10-05 14:33:18.375823"
    5  = "2021-10-05 14:33:18.375826"
    6  = "2021-10-05 14:33:18.375829"
    7  = "2021-10-05 14:33:18.375831"
    8  = "2021-10-05 14:33:18.375834"
    9  = "2021-10-05 14:33:18.375836"
    10 = "2021-10-05 14:33:18.375839"
    11 = "2021-10-05 14:33:18.375841"
    12 = "2021-10-05 14:33:18.375844"
    13 = "2021-10-05 14:33:18.375846"
    14 = "2021-10-05 14:33:18.375849"
    15 = "2021-10-05 14:33:18.375851"
    16 = "2021-10-05 14:33:18.375854"
    17 = "2021-10-05 14:33:18.375856"
    18 = "2021-10-05 14:33:18.375859"
    19 = "2021-10-05 14:33:18.375862"
    20 = "2021-10-05 14:33:18.375865"
    21 = "2021-10-05 14:33:18.375867"
    22 = "2021-10-05 14:33:18.375870"
    23 = "2021-10-05 14:33:18.375872"
}

foreach ($rowNum in $timestamps.Keys) {
    $dataSheet.Range("F$rowNum").Value = $timestamps[$rowNum]
}

# --- Add a new "metadata" worksheet right after the "data" sheet ---
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (bold/centered/bordered style copied from the data sheet header)
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

$dataSheet.Range("B1").Copy()
$metaSheet.Range("B1:G1").PasteSpecial(-4122)

# Data row
$metaSheet.Range("B2").Value = "Brain Channelopathies"
$metaSheet.Range("C2").Value = 74
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "1.0"
$metaSheet.Range("E2").Value = "2020-12-30T01:37:05.618866Z"
$metaSheet.Range("F2").NumberFormat = "@"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:18.372532"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/74/?format=json"

$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$metaSheet.Range("A2").Value = 0

$dataSheet.Activate()
$dataSheet.Range("A1").Select()

Write-Output "metadata sheet added; timestamps updated"
